$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("106:110").Delete()
